$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating point drift on the existing last row (58)
$ws.Range("A58").Value = 44371.76033399074

# Append new row 59 with the latest retrieved data
$ws.Range("A59").Value = 44372.76793294882
$ws.Range("A59").NumberFormat = $ws.Range("A58").NumberFormat
$ws.Range("B59").Value = 79164
$ws.Range("C59").Value = 66531
$ws.Range("D59").Value = 3584
$ws.Range("E59").Value = 2158
$ws.Range("F59").Value = 1539
$ws.Range("G59").Value = 20944
$ws.Range("H59").Value = 1491
$ws.Range("I59").Value = 895
$ws.Range("J59").Value = 186
